# Update LR-pairs sheet with newly computed TPM-based NATMI values.
# Sending clusters ECs/FAPs/MuSCs x 3 targets (9 rows, rows 2-10) collapses
# to sending clusters FAPs/MuSCs x 3 targets (6 rows, rows 2-7); the "ECs"
# sending-cluster rows are dropped and every remaining row gets refreshed
# numeric values from the new TPM run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for rows 2-7 (columns A..T), keyed by row number.
$rowValues = @{
    2 = @("FAPs","Sema3c","Plxnd1","ECs",
          3,1,39.09670133333334,117.290104,0.9758026125363394,0.9758026125363395,
          3,1,53.62555,160.87665,0.5942801874957572,0.5942801874957572,
          2096.582112185733,18869.2390096716,0.5799001595369454,0.5799001595369455)
    3 = @("FAPs","Sema3c","Plxnd1","FAPs",
          3,1,39.09670133333334,117.290104,0.9758026125363394,0.9758026125363395,
          3,1,8.866675333333333,26.600026,0.09826080067350991,0.0982608006735099,
          346.6577573269671,3119.919815942704,0.09588314600712347,0.09588314600712346)
    4 = @("FAPs","Sema3c","Plxnd1","MuSCs",
          3,1,39.09670133333334,117.290104,0.9758026125363394,0.9758026125363395,
          3,1,27.74391433333333,83.23174299999999,0.3074590118307329,0.3074590118307329,
          1084.695532507919,9762.259792571273,0.3000193069922705,0.3000193069922705)
    5 = @("MuSCs","Sema3c","Plxnd1","ECs",
          3,1,0.9694973333333333,2.908492,0.02419738746366056,0.02419738746366056,
          3,1,53.62555,160.87665,0.5942801874957572,0.5942801874957572,
          51.98982772353333,467.9084495117999,0.01438002795881168,0.01438002795881168)
    6 = @("MuSCs","Sema3c","Plxnd1","FAPs",
          3,1,0.9694973333333333,2.908492,0.02419738746366056,0.02419738746366056,
          3,1,8.866675333333333,26.600026,0.09826080067350991,0.0982608006735099,
          8.596218091199111,77.36596282079199,0.002377654666386437,0.002377654666386437)
    7 = @("MuSCs","Sema3c","Plxnd1","MuSCs",
          3,1,0.9694973333333333,2.908492,0.02419738746366056,0.02419738746366056,
          3,1,27.74391433333333,83.23174299999999,0.3074590118307329,0.3074590118307329,
          26.89765096239511,242.078858661556,0.007439704838462439,0.007439704838462438)
}

foreach ($r in 2..7) {
    $vals = $rowValues[$r]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Rows 8-10 (old MuSCs-sending data, now folded into rows 5-7 above) are
# no longer needed - remove them so the used range shrinks back to T7.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
